$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update B3 value; A3 formula (=B3/B1*A1) will recalculate automatically.
$ws.Range("B3").Value = 1100
